# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header "K") needs to be rewritten with freshly
# computed strikeout values (s_vals) for each game row. Rows that are not
# present in the new s_vals map are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values keyed by worksheet row number.
$s_vals = [ordered]@{
    2  = 1
    3  = 1
    4  = 1
    5  = 2
    6  = 2
    7  = 0
    8  = 1
    9  = 3
    10 = 1
    12 = 1
    13 = 1
    14 = 2
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 2
    20 = 2
    21 = 1
    22 = 0
    23 = 0
    24 = 2
    25 = 0
    26 = 0
    27 = 1
    28 = 2
    29 = 2
    30 = 1
    31 = 0
    32 = 0
    33 = 2
    34 = 2
    35 = 2
    36 = 1
    37 = 2
    38 = 2
    39 = 1
    40 = 2
    41 = 1
    42 = 1
    43 = 3
    44 = 2
    45 = 1
    46 = 0
    47 = 1
    48 = 0
    49 = 1
    50 = 1
    51 = 1
    52 = 1
    53 = 3
    54 = 2
    55 = 2
    56 = 0
    57 = 1
    58 = 1
    59 = 1
    60 = 1
    61 = 1
    62 = 2
    63 = 1
    64 = 3
    65 = 0
    66 = 2
    67 = 0
    69 = 1
    70 = 1
    71 = 1
    72 = 1
    73 = 0
    74 = 1
    75 = 0
    76 = 0
    77 = 3
    78 = 2
    79 = 3
    82 = 1
    83 = 2
    84 = 0
    85 = 0
    86 = 1
}

foreach ($row in $s_vals.Keys) {
    $ws.Cells.Item($row, 7).Value = $s_vals[$row]
}
